$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string labels / small data points on the existing table ---
# (order chosen to reproduce the sharedStrings.xml insertion order: "po tuningu:",
#  "sr zapytania", "sr po tuningu", "po jdbc.batch")
$ws.Range("D59").Value = "po tuningu:"
$ws.Range("D56").Value = "śr zapytania"
$ws.Range("D65").Value = "śr po tuningu"
$ws.Range("F18").Value = "po jdbc.batch"
$ws.Range("F19").Value = 48917

# --- Helper table used as the source for the new "po tuningu" bar chart ---
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 2
$ws.Range("F60").Value = 3
$ws.Range("G60").Value = 4
$ws.Range("H60").Value = 5
$ws.Range("I60").Value = 6

$ws.Range("D61").Value = 280
$ws.Range("E61").Value = 310
$ws.Range("F61").Value = 490
$ws.Range("G61").Value = 480
$ws.Range("H61").Value = 690
$ws.Range("I61").Value = 610

$ws.Range("D62").Value = 210
$ws.Range("E62").Value = 300
$ws.Range("F62").Value = 330
$ws.Range("G62").Value = 550
$ws.Range("H62").Value = 680
$ws.Range("I62").Value = 700

$ws.Range("D63").Value = 250
$ws.Range("E63").Value = 230
$ws.Range("F63").Value = 260
$ws.Range("G63").Value = 420
$ws.Range("H63").Value = 680
$ws.Range("I63").Value = 530

$ws.Range("D64").Value = 250
$ws.Range("E64").Value = 230
$ws.Range("F64").Value = 310
$ws.Range("G64").Value = 680
$ws.Range("H64").Value = 520
$ws.Range("I64").Value = 580

$ws.Range("D66").Value = 247
$ws.Range("E66").Value = 267
$ws.Range("F66").Value = 347
$ws.Range("G66").Value = 532
$ws.Range("H66").Value = 642
$ws.Range("I66").Value = 605

# --- Column D got a touch wider once the new label text was added ---
$ws.Columns.Item(4).ColumnWidth = 12

# --- Reposition / resize the second ("po tuningu") line chart to make room ---
$co2 = $ws.ChartObjects().Item(2)
$co2.Left = 212.17878967765748
$co2.Top = 330
$co2.Width = 497.1288275098425
$co2.Height = 182.25

# --- Add the new clustered-column chart comparing "normalnie" vs "tuning" ---
$co4 = $ws.ChartObjects().Add(250.3662109375, 513.75, 372.697265625, 216)
$chart4 = $co4.Chart
$chart4.ChartType = 51
$chart4.SeriesCollection().NewSeries()
$s1 = $chart4.SeriesCollection(1)
$s1.Name = "normalnie"
$s1.Values = $ws.Range("D57:I57")
$chart4.SeriesCollection().NewSeries()
$s2 = $chart4.SeriesCollection(2)
$s2.Name = "tuning"
$s2.Values = $ws.Range("D66:I66")
$chart4.HasLegend = $true
$chart4.Legend.Position = -4152

$co4.Name = "Wykres 4"
$co4.Top = 513.75
$co4.Left = 250.3662109375
$co4.Width = 372.697265625
$co4.Height = 216

# --- Move the selection/view to roughly where the author left it ---
$ws.Range("G20").Select()

Write-Host "edit applied"
